# Update KM gas composition to include averaged values and standard deviation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 16 (refill #15, the Thanksgiving refill) gets its previously-missing
#     dates filled in (11/26/2022), and its note updated to describe the fix ---
$ws.Range("B2").Copy()
$ws.Range("B16").PasteSpecial(-4122)
$ws.Range("B16").Value = 44891

$ws.Range("G2").Copy()
$ws.Range("G16").PasteSpecial(-4122)
$ws.Range("G16").Value = 44891

$ws.Range("E16").Value = "Thanksgiving refill, average 11/26 and 11/27"

# --- Rename / re-purpose the header row (row 1) to snake_case field names ---
$ws.Range("A1").Value = "refill_no"
$ws.Range("B1").Value = "rawhide_refill_date"
$ws.Range("C1").Value = "start_utc"
$ws.Range("D1").Value = "end_utc"
$ws.Range("E1").Value = "notes"
$ws.Range("F1").Value = "su_notes"
$ws.Range("G1").Value = "km_composition_date"
# H1 (su_raw) and I1 (su_normalized) stay as-is.

# The su_km column (J) header & data are no longer populated - clear it out entirely.
$ws.Range("J1").ClearContents()
$ws.Range("J2:J16").ClearContents()

# --- Restore the last active-cell selection to match the editor's cursor ---
[void]$ws.Range("J12").Select()
